# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to match the freshly generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 5498
$ws1.Range("F23").Value = 3551
$ws1.Range("F24").Value = 1120
$ws1.Range("F25").Value = 2816
$ws1.Range("F27").Value = 1981
$ws1.Range("F28").Value = 4074
$ws1.Range("F30").Value = 917
$ws1.Range("F32").Value = 1297
$ws1.Range("F33").Value = 64
$ws1.Range("F36").Value = 1270
$ws1.Range("F38").Value = 1047
$ws1.Range("F39").Value = 676
$ws1.Range("F42").Value = 34
$ws1.Range("F43").Value = 313

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5498
$ws4.Range("F23").Value = 3551
$ws4.Range("F26").Value = 1120
$ws4.Range("F27").Value = 2816
$ws4.Range("F28").Value = 1981
$ws4.Range("F29").Value = 4074
$ws4.Range("F32").Value = 917
$ws4.Range("F33").Value = 1297
$ws4.Range("F37").Value = 1270
$ws4.Range("F39").Value = 1047
$ws4.Range("F41").Value = 676
$ws4.Range("F47").Value = 313
